$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1

# Add new row 3 with values, copying A2's style (bold font, border, centered) to A3
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 1
